$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RF (column I) values for rows 25 through 44 to the new value
$ws.Range("I25:I44").Value = 7.152777777777778
